$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 360.91666
$ws.Range("I28").Value = 91.2
$ws.Range("J28").Value = 1709.5
$ws.Range("K28").Value = 91.2
$ws.Range("L28").Value = 1709.5
$ws.Range("M28").Value = 393.8
$ws.Range("N28").Value = -2679.5
$ws.Range("H87").Value = 59968.25
$ws.Range("J87").Value = 99949
$ws.Range("L87").Value = 99949
$ws.Range("N87").Value = -102445
$ws.Range("H90").Value = 59968.25
$ws.Range("J90").Value = 99949
$ws.Range("L90").Value = 299847
$ws.Range("N90").Value = -312327
$ws.Range("H94").Value = 11000.5
$ws.Range("I94").Value = 11000.5
$ws.Range("K94").Value = 11000.5
$ws.Range("M94").Value = -10549.5
$ws.Range("H98").Value = 964.5
$ws.Range("I98").Value = 857.8
$ws.Range("K98").Value = 857.8
$ws.Range("M98").Value = 640.2
$ws.Range("H122").Value = 964.5
$ws.Range("I122").Value = 857.8
$ws.Range("K122").Value = 2573.4
$ws.Range("M122").Value = -123.3999999999996
$ws.Range("H137").Value = 487.25
$ws.Range("J137").Value = 1000
$ws.Range("L137").Value = 3000
$ws.Range("N137").Value = -8100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1455
$ws.Range("I2").Value = 1455
$ws.Range("K2").Value = 1455
$ws.Range("M2").Value = -1342
$ws.Range("H32").Value = 5271
$ws.Range("I32").Value = 5271
$ws.Range("K32").Value = 5271
$ws.Range("M32").Value = -4984
$ws.Range("H102").Value = 1942.4286
$ws.Range("I102").Value = 1399.4
$ws.Range("J102").Value = 3300
$ws.Range("K102").Value = 1399.4
$ws.Range("L102").Value = 3300
$ws.Range("M102").Value = 222.5999999999999
$ws.Range("N102").Value = -6544
$ws.Range("H110").Value = 1997.5
$ws.Range("I110").Value = 995.5
$ws.Range("J110").Value = 2999.5
$ws.Range("K110").Value = 995.5
$ws.Range("L110").Value = 2999.5
$ws.Range("M110").Value = 1049.5
$ws.Range("N110").Value = -7089.5
$ws.Range("H116").Value = 1455
$ws.Range("I116").Value = 1455
$ws.Range("K116").Value = 1455
$ws.Range("M116").Value = 839
$ws.Range("H122").Value = 2933.3333
$ws.Range("I122").Value = 1800
$ws.Range("K122").Value = 5400
$ws.Range("M122").Value = -2950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1455
$ws.Range("I3").Value = 1455
$ws.Range("K3").Value = 1455
$ws.Range("M3").Value = -1341
$ws.Range("H86").Value = 2045.25
$ws.Range("I86").Value = 1884.7142
$ws.Range("K86").Value = 1884.7142
$ws.Range("M86").Value = -761.7141999999999
$ws.Range("H89").Value = 2045.25
$ws.Range("I89").Value = 1884.7142
$ws.Range("K89").Value = 9423.571
$ws.Range("M89").Value = -3807.571
$ws.Range("H94").Value = 2117.125
$ws.Range("I94").Value = 1479.1666
$ws.Range("J94").Value = 2499.9
$ws.Range("K94").Value = 1479.1666
$ws.Range("L94").Value = 2499.9
$ws.Range("M94").Value = -1028.1666
$ws.Range("N94").Value = -3401.9
$ws.Range("H99").Value = 3229.6667
$ws.Range("I99").Value = 3229.6667
$ws.Range("K99").Value = 3229.6667
$ws.Range("M99").Value = -1731.6667
$ws.Range("H102").Value = 7733.1665
$ws.Range("I102").Value = 7733.1665
$ws.Range("K102").Value = 7733.1665
$ws.Range("M102").Value = -4488.1665
$ws.Range("H105").Value = 2390.7273
$ws.Range("I105").Value = 2390.7273
$ws.Range("K105").Value = 2390.7273
$ws.Range("M105").Value = -643.7273
$ws.Range("H111").Value = 45000
$ws.Range("J111").Value = 45000
$ws.Range("L111").Value = 45000
$ws.Range("N111").Value = -53180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1076.8462
$ws.Range("I16").Value = 910.36365
$ws.Range("K16").Value = 910.36365
$ws.Range("M16").Value = -623.36365
$ws.Range("H113").Value = 1076.8462
$ws.Range("I113").Value = 910.36365
$ws.Range("K113").Value = 910.36365
$ws.Range("M113").Value = 1259.63635
$ws.Range("H122").Value = 4734.7144
$ws.Range("I122").Value = 764.6667
$ws.Range("K122").Value = 2294.0001
$ws.Range("M122").Value = 155.9998999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3168.1667
$ws.Range("I80").Value = 3249.5
$ws.Range("K80").Value = 3249.5
$ws.Range("M80").Value = -2251.5
$ws.Range("H83").Value = 3168.1667
$ws.Range("I83").Value = 3249.5
$ws.Range("K83").Value = 16247.5
$ws.Range("M83").Value = -11255.5
$ws.Range("H97").Value = 1283.3334
$ws.Range("I97").Value = 1300
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 1300
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -804
$ws.Range("N97").Value = -2242
$ws.Range("H132").Value = 1471.6666
$ws.Range("I132").Value = 1471.6666
$ws.Range("K132").Value = 4414.9998
$ws.Range("M132").Value = -1884.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2938.3572
$ws.Range("I82").Value = 1769.75
$ws.Range("J82").Value = 4496.5
$ws.Range("K82").Value = 1769.75
$ws.Range("L82").Value = 4496.5
$ws.Range("M82").Value = -1408.75
$ws.Range("N82").Value = -5218.5
$ws.Range("H85").Value = 2938.3572
$ws.Range("I85").Value = 1769.75
$ws.Range("J85").Value = 4496.5
$ws.Range("K85").Value = 1769.75
$ws.Range("L85").Value = 4496.5
$ws.Range("M85").Value = -521.75
$ws.Range("N85").Value = -6992.5
$ws.Range("H93").Value = 2436.2703
$ws.Range("I93").Value = 2436.2703
$ws.Range("K93").Value = 2436.2703
$ws.Range("M93").Value = -1188.2703
$ws.Range("H100").Value = 4450
$ws.Range("I100").Value = 4500
$ws.Range("K100").Value = 4500
$ws.Range("M100").Value = -3959

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3277.4443
$ws.Range("I81").Value = 3374.75
$ws.Range("J81").Value = 2499
$ws.Range("K81").Value = 6749.5
$ws.Range("L81").Value = 4998
$ws.Range("M81").Value = -5688.5
$ws.Range("N81").Value = -7120
$ws.Range("H84").Value = 3277.4443
$ws.Range("I84").Value = 3374.75
$ws.Range("J84").Value = 2499
$ws.Range("K84").Value = 33747.5
$ws.Range("L84").Value = 24990
$ws.Range("M84").Value = -28443.5
$ws.Range("N84").Value = -35598
$ws.Range("H100").Value = 906.625
$ws.Range("I100").Value = 875.5
$ws.Range("K100").Value = 1751
$ws.Range("M100").Value = -1210
$ws.Range("H107").Value = 616.6667
$ws.Range("H113").Value = 331.66666
$ws.Range("I113").Value = 347.5
$ws.Range("J113").Value = 300
$ws.Range("K113").Value = 1042.5
$ws.Range("L113").Value = 900
$ws.Range("M113").Value = 1127.5
$ws.Range("N113").Value = -5240
$ws.Range("H122").Value = 1184.625
$ws.Range("I122").Value = 895.8
$ws.Range("K122").Value = 2687.4
$ws.Range("M122").Value = -237.3999999999996
$ws.Range("H132").Value = 2642.4285
$ws.Range("I132").Value = 2582.8333
$ws.Range("K132").Value = 7748.499899999999
$ws.Range("M132").Value = -5218.499899999999
